$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K (11) and add the "categories" header.
$null = $ws.Columns("K:K").Insert()
$ws.Range("K1").Value2 = "categories"

# The AutoFilter range does not auto-expand when a column is inserted in the
# middle of it, so re-apply it across the new full range, preserving the
# existing discrete-value filter on column B (colId 1 -> field 2).
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:AH56").AutoFilter(2, @("C53630"), 7)

# Keep the workbook-level _FilterDatabase defined name in sync with the
# expanded AutoFilter range.
$fdb = $wb.Names.Item("Collection_CM!_FilterDatabase")
$fdb.RefersTo = "=Collection_CM!`$A`$1:`$AH`$56"

# Restore the previous selection (shifted one column right because of the
# inserted column).
$null = $ws.Range("K36").Select()
